$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sensor")

$ws.Range("D11").Value = 4
$ws.Range("D12").Value = 2

$ws2 = $wb.Worksheets.Item("HAT")
$ws2.Range("G8").Select()
$ws.Activate()
